# New weekly price-report record for "Vega Modelo de Temuco - Albahaca".
# A new observation (Fecha 44754, Volumen 30) is inserted as row 193,
# pushing the existing rows 193:245 down to 194:246.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 193:245 down one row, duplicating row 193's formatting
# (matches the existing D-column date style) into the freshly opened row.
$ws.Rows.Item(193).Insert()

# Populate the new row 193 with the new weekly record. Columns not called
# out below (K, L, M, N, O, P, Q) keep the same values the prior row 193
# already had (6000 / 6000 / 6000 / $/paquete / Región de Arica y
# Parinacota / 6000 / 1), so only the date and volume actually change.
$ws.Range("A193").Value = 10
$ws.Range("B193").Value = "Vega Modelo de Temuco"
$ws.Range("C193").Value = "La Araucanía"
$ws.Range("D193").Value = 44754
$ws.Range("E193").Value = 9
$ws.Range("F193").Value = 100112052
$ws.Range("G193").Value = "Albahaca"
$ws.Range("H193").Value = "Sin especificar"
$ws.Range("I193").Value = "Primera"
$ws.Range("J193").Value = 30
$ws.Range("K193").Value = 6000
$ws.Range("L193").Value = 6000
$ws.Range("M193").Value = 6000
$ws.Range("N193").Value = "$/paquete"
$ws.Range("O193").Value = "Región de Arica y Parinacota"
$ws.Range("P193").Value = 6000
$ws.Range("Q193").Value = 1
$ws.Range("R193").Value = "Hortaliza"
